# Fill in a new item row (row 4) on the pharmacy transactions sheet, mark the
# "Name" and "Current balance" merged cells as Text-formatted so values like
# "2:0" are preserved verbatim, carry the price total down into the footer
# (K5), and grow row 5 to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Name (B4:G4) and Current balance (H4:K4) merged ranges as Text
# before typing into them, so Excel doesn't reinterpret values such as "2:0".
$ws.Range("B4:G4").NumberFormat = "@"
$ws.Range("H4:K4").NumberFormat = "@"

# Row 4 - new transaction entry
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "RHINEX 0.05% INFANTILE NASAL DROPS 10 ML"
$ws.Range("H4").Value = "2:0"
$ws.Range("L4").Value = 18
$ws.Range("N4").Value = 1

# Row 5 - totals footer
$ws.Range("K5").Value = 18
$ws.Rows.Item(5).RowHeight = 26.25
